# Apply the WIP tracker update:
#  - Consolidate the "Campaign (group)" column (A3:A12) so every COP row
#    rolls up under a single "COP 2026" label (was a mix of
#    COP 2025 / COP 2026 / COP 2027 / COP 2028).
#  - Lengthen the Digital Display Phase 1 details note (E5).
#  - Move/extend the active selection to A3:A12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll every campaign-group cell in rows 3-12 up to "COP 2026".
$ws.Range("A3:A12").Value = "COP 2026"

# Expand the details text on row 5 (Digital Display Phase 1).
$ws.Range("E5").Value = "50% approved, 50% in review, working on Phase 2, this is also test how much the text is flowing. Normally this wouldn't be this long"

# Update the selection shown when the sheet re-opens.
$ws.Range("A3:A12").Select()
